# Rename the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Sheet"

# Clear the old layout (column D / rows beyond the new data / row 10 placeholder)
$ws.Cells.Clear()

# Header row
$ws.Range("A1").Value = "번호"
$ws.Range("B1").Value = "영어"
$ws.Range("C1").Value = "수학"

# Data rows (번호, 영어, 수학)
$data = @(
    @(1, 77, 44),
    @(2, 98, 74),
    @(3, 14, 31),
    @(4, 100, 55),
    @(5, 46, 75),
    @(6, 14, 82),
    @(7, 11, 32),
    @(8, 89, 43),
    @(9, 7, 39),
    @(10, 63, 59)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
